$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while forcing Excel to treat it as literal
# text (matching the source file's inline-string cells), instead of letting
# the normal Range.Value setter auto-detect numeric-looking strings (e.g.
# "69.231.29", "2.16") and convert them into Number cells. We do this by
# building the text via a formula ("="<value>"") in a scratch cell, copying
# it, and pasting only the resulting value into the destination - this keeps
# the destination cell's string-ness without leaving any NumberFormat /
# quote-prefix style residue behind.
function Set-TextValue {
    param($ws, $addr, $val)
    $scratch = $ws.Range("ZZ1")
    $scratch.Formula = '="' + $val + '"'
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)   # xlPasteValues
    $scratch.Value = $null
}

Set-TextValue $ws 'D2' '69.231.29'
$ws.Range('E2').Value = '  -3.00%  '
Set-TextValue $ws 'D3' '2.487.43'
$ws.Range('E3').Value = '  -3.44%  '
$ws.Range('E4').Value = '  -0.06%  '
Set-TextValue $ws 'D5' '565.86'
$ws.Range('E5').Value = '  -3.17%  '
Set-TextValue $ws 'D6' '164.14'
$ws.Range('E6').Value = '  -5.49%  '
$ws.Range('E7').Value = '  -0.09%  '
Set-TextValue $ws 'D8' '0.511'
$ws.Range('E8').Value = '  -1.80%  '
Set-TextValue $ws 'D9' '2.486.05'
$ws.Range('E9').Value = '  -3.45%  '
Set-TextValue $ws 'D10' '0.158'
$ws.Range('E10').Value = '  -5.85%  '
$ws.Range('E11').Value = '  -0.88%  '
Set-TextValue $ws 'D12' '0.353'
$ws.Range('E12').Value = '  -1.78%  '
Set-TextValue $ws 'D13' '4.90'
$ws.Range('E13').Value = '  -0.49%  '
Set-TextValue $ws 'D14' '2.939.25'
$ws.Range('E14').Value = '  -3.69%  '
Set-TextValue $ws 'D15' '69.107.30'
$ws.Range('E15').Value = '  -3.02%  '
Set-TextValue $ws 'D16' '0.0000175'
$ws.Range('E16').Value = '  -3.67%  '
Set-TextValue $ws 'D17' '24.22'
$ws.Range('E17').Value = '  -5.54%  '
Set-TextValue $ws 'D18' '2.491.86'
$ws.Range('E18').Value = '  -3.47%  '
Set-TextValue $ws 'D19' '11.12'
$ws.Range('E19').Value = '  -4.70%  '
Set-TextValue $ws 'D20' '7.37'
Set-TextValue $ws 'D21' '345.53'
$ws.Range('E21').Value = '  -3.81%  '
Set-TextValue $ws 'D22' '3.86'
$ws.Range('E22').Value = '  -3.14%  '
$ws.Range('E23').Value = '  -7.91%  '
$ws.Range('E24').Value = '  +0.01%  '
Set-TextValue $ws 'D25' '69.39'
$ws.Range('E25').Value = '  -1.75%  '
$ws.Range('E26').Value = '  -6.50%  '
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws 'D27' '8.66'
$ws.Range('E27').Value = '  -6.18%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
Set-TextValue $ws 'D28' '2.609.43'
$ws.Range('E28').Value = '  -2.87%  '
Set-TextValue $ws 'D29' '0.999'
$ws.Range('E29').Value = '  -0.42%  '
Set-TextValue $ws 'D30' '0.0₃0872'
$ws.Range('E30').Value = '  -6.45%  '
Set-TextValue $ws 'D31' '7.70'
$ws.Range('E31').Value = '  -4.00%  '
Set-TextValue $ws 'D32' '442.86'
$ws.Range('E32').Value = '  -7.44%  '
Set-TextValue $ws 'D33' '1.19'
$ws.Range('E33').Value = '  -7.90%  '
$ws.Range('E34').Value = '  -0.02%  '
$ws.Range('E35').Value = '  -4.83%  '
Set-TextValue $ws 'D36' '155.02'
$ws.Range('E36').Value = '  -2.04%  '
$ws.Range('E37').Value = '  -4.78%  '
$ws.Range('E38').Value = '  -0.54%  '
Set-TextValue $ws 'D39' '18.11'
$ws.Range('E39').Value = '  -4.15%  '
$ws.Range('E40').Value = '  +0.01%  '
Set-TextValue $ws 'D41' '0.313'
$ws.Range('E41').Value = '  -3.48%  '
Set-TextValue $ws 'D42' '4.57'
$ws.Range('E42').Value = '  -7.40%  '
$ws.Range('E43').Value = '  -4.21%  '
Set-TextValue $ws 'D44' '37.90'
$ws.Range('E44').Value = '  -2.07%  '
Set-TextValue $ws 'D45' '2.16'
Set-TextValue $ws 'D46' '1.07'
$ws.Range('E46').Value = '  -10.07%  '
Set-TextValue $ws 'D47' '139.61'
$ws.Range('E47').Value = '  -4.75%  '
Set-TextValue $ws 'D48' '3.43'
$ws.Range('E48').Value = '  -4.02%  '
Set-TextValue $ws 'D49' '0.512'
$ws.Range('E49').Value = '  -5.80%  '
Set-TextValue $ws 'D50' '0.0723'
$ws.Range('E50').Value = '  -2.60%  '
Set-TextValue $ws 'D51' '0.572'
$ws.Range('E51').Value = '  -3.18%  '

$excel.CutCopyMode = $false
